# Workbook edit: add Sheet2, update Sheet1!A2 timestamp, make Sheet2 the active sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update the timestamp value in Sheet1 (A2), preserving its existing date/time format.
$ws1.Range("A2").Value = 43002.6080540257

# Insert a new worksheet named "Sheet2" right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws2.Range("A1").Value = 20

# Make the newly added sheet the active / selected tab.
$ws2.Activate()
